# This script reorders rows 2-15 (columns C:F -> runs/balls/fours/sixes)
# on the "Steven Smith (c)" sheet to match the updated per-innings ordering,
# while columns A (playerName) and B (teamName) are identical across rows
# and therefore unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage 1: snapshot current C:F values for rows 2-15 into staging rows (offset +100)
# so the permutation below can be applied without clobbering source data.
$ws.Range("C2:F2").Copy($ws.Range("C102:F102"))
$ws.Range("C3:F3").Copy($ws.Range("C103:F103"))
$ws.Range("C4:F4").Copy($ws.Range("C104:F104"))
$ws.Range("C5:F5").Copy($ws.Range("C105:F105"))
$ws.Range("C6:F6").Copy($ws.Range("C106:F106"))
$ws.Range("C7:F7").Copy($ws.Range("C107:F107"))
$ws.Range("C8:F8").Copy($ws.Range("C108:F108"))
$ws.Range("C9:F9").Copy($ws.Range("C109:F109"))
$ws.Range("C10:F10").Copy($ws.Range("C110:F110"))
$ws.Range("C11:F11").Copy($ws.Range("C111:F111"))
$ws.Range("C12:F12").Copy($ws.Range("C112:F112"))
$ws.Range("C13:F13").Copy($ws.Range("C113:F113"))
$ws.Range("C14:F14").Copy($ws.Range("C114:F114"))
$ws.Range("C15:F15").Copy($ws.Range("C115:F115"))

# Stage 2: copy staged rows into their final destination row per the new order.
$ws.Range("C111:F111").Copy($ws.Range("C2:F2"))
$ws.Range("C112:F112").Copy($ws.Range("C3:F3"))
$ws.Range("C102:F102").Copy($ws.Range("C4:F4"))
$ws.Range("C103:F103").Copy($ws.Range("C5:F5"))
$ws.Range("C104:F104").Copy($ws.Range("C6:F6"))
$ws.Range("C109:F109").Copy($ws.Range("C7:F7"))
$ws.Range("C113:F113").Copy($ws.Range("C8:F8"))
$ws.Range("C110:F110").Copy($ws.Range("C9:F9"))
$ws.Range("C106:F106").Copy($ws.Range("C10:F10"))
$ws.Range("C107:F107").Copy($ws.Range("C11:F11"))
$ws.Range("C115:F115").Copy($ws.Range("C12:F12"))
$ws.Range("C105:F105").Copy($ws.Range("C13:F13"))
$ws.Range("C114:F114").Copy($ws.Range("C14:F14"))
$ws.Range("C108:F108").Copy($ws.Range("C15:F15"))

# Stage 3: clear the temporary staging rows.
$ws.Range("C102:F115").ClearContents()
